# This script re-orders the data rows (2..22) of the active sheet.
# The edit described by the diff is a pure row permutation: each target
# row receives the full set of column D..R values that used to live in
# another (source) row. Rows 12, 13 and 14 are unchanged (map to themselves).
#
# Mapping is new_row -> source_row (i.e. "new row N should contain what
# used to be in source_row before the edit").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$perm = @{}
$perm[2]  = 20
$perm[3]  = 9
$perm[4]  = 10
$perm[5]  = 21
$perm[6]  = 16
$perm[7]  = 17
$perm[8]  = 18
$perm[9]  = 19
$perm[10] = 11
$perm[11] = 2
$perm[12] = 12
$perm[13] = 13
$perm[14] = 14
$perm[15] = 4
$perm[16] = 3
$perm[17] = 6
$perm[18] = 7
$perm[19] = 8
$perm[20] = 22
$perm[21] = 5
$perm[22] = 15

# Columns D (4) through R (18) hold the per-row record data that moves
# together as a unit (A, B, C, Q, R are identical across all rows so they
# do not matter, but we keep them in the captured range for safety).
$firstCol = 1
$lastCol  = 18

# 1) Snapshot every row's current values before any writes happen, since
#    this is an in-place permutation (source rows get overwritten too).
$snapshot = @{}
for ($row = 2; $row -le 22; $row++) {
    $rowData = @{}
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $rowData[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshot[$row] = $rowData
}

# 2) Write back each target row using the snapshot of its source row.
for ($row = 2; $row -le 22; $row++) {
    $srcRow = $perm[$row]
    $rowData = $snapshot[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col]
    }
}
